$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting existing rows 48:96 down to 49:97
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new record
$ws.Cells.Item(48, 1).Value = 7
$ws.Cells.Item(48, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(48, 3).Value = "Ñuble"
$ws.Cells.Item(48, 4).Value = 45240
$ws.Cells.Item(48, 5).Value = 16
$ws.Cells.Item(48, 6).Value = 100112026
$ws.Cells.Item(48, 7).Value = "Haba"
$ws.Cells.Item(48, 8).Value = "Sin especificar"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 60
$ws.Cells.Item(48, 11).Value = 10000
$ws.Cells.Item(48, 12).Value = 10000
$ws.Cells.Item(48, 13).Value = 10000
$ws.Cells.Item(48, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(48, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(48, 16).Value = 400
$ws.Cells.Item(48, 17).Value = 25
$ws.Cells.Item(48, 18).Value = "Hortaliza"
